$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20: Exportaciones 0%, pulled from L12 (Exportaciones column), plain currency format.
$ws.Range("L20").Value = "Exportaciones 0%"
$ws.Range("M20").Formula = "=L12"
$ws.Range("M20").NumberFormat = "_-""$""* #,##0.00_-;\-""$""* #,##0.00_-;_-""$""* ""-""??_-;_-@_-"

# Row 16: label changes to "Ventas Gravadas 13%:" and its total now pulls straight
# from K12 (Ventas Gravadas) instead of SUM(K12:L12) (which used to fold in Exportaciones).
$ws.Range("L16").Value = "Ventas Gravadas 13%:"
$ws.Range("M16").Formula = "=K12"

# New row 21: TOTAL, bold label + bold currency total combining all the lines above.
$ws.Range("L21").Value = "TOTAL"
$ws.Range("L21").Font.Bold = $true
$ws.Range("M21").Formula = "=M14+M15+M16-M19+M20"
$ws.Range("M21").Font.Bold = $true
$ws.Range("M21").NumberFormat = "_-""$""* #,##0.00_-;\-""$""* #,##0.00_-;_-""$""* ""-""??_-;_-@_-"

# Selection state left on L22 to match the source session.
$null = $ws.Range("L22").Select()
